$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: turn the placeholder/blank row into a new "PPO2" / FULL data row ---
$ws.Range("A11").Value = "PPO2"
$ws.Range("B11").Value = "FULL"
$ws.Range("C11").Value = 14120.883959999999
$ws.Range("D11").Value = 21060.503959999998
$ws.Range("E11").Value = 1532.224467
$ws.Range("F11").Value = 17096.225419999999
$ws.Range("F11").Font.Bold = $false
# H11 already carries the shared "(F#-10000)/10000 * (365/K#)" formula and
# recalculates on its own once F11/K3 are known - leave it alone.

# Re-stamp the percent format on the whole Yearly-Return column so it lands on
# a fresh style record (matches the saved workbook, which reshuffled cellXfs
# after the font/numFmt table was touched) without altering the visible format.
$ws.Range("H2:H11").Orientation = 0

# --- Row 12: drop the old "yearly return delta" formula entirely (cell + format) ---
# it gets re-homed to row 16 below.
$ws.Range("H12").Clear()

# --- Row 16: re-home the final MIN-strategy comparison row ---
$ws.Range("F16").Value = 15176.64
$ws.Range("H16").Formula = "=H3-H11"
$ws.Range("H16").NumberFormat = "0.0%"

# --- Selection bookkeeping (matches the saved file) ---
$ws.Range("D12").Select()
